$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition): update "想去人数" (F column) counts for rows 2-4
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 984
$wsExpo.Range("F3").Value = 1997
$wsExpo.Range("F4").Value = 443

# Sheet "全部类型" (All types): update the same three events (rows 4-6)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 984
$wsAll.Range("F5").Value = 1997
$wsAll.Range("F6").Value = 443
